$d = $word.ActiveDocument

# Add solid black borders (outside + inside, 1/2 pt = sz 4) to every table
# in the document, matching a <w:tblBorders> block of:
#   <w:top/left/bottom/right/insideH/insideV w:val="single" w:sz="4" w:space="0" w:color="000000"/>
foreach ($t in $d.Tables) {
    $t.Borders.Enable = $true
    $t.Borders.OutsideLineStyle = 1   # wdLineStyleSingle
    $t.Borders.OutsideLineWidth = 2   # wdLineWidth025pt -> w:sz="4"
    $t.Borders.OutsideColor = 0       # wdColorBlack -> w:color="000000"
    $t.Borders.InsideLineStyle = 1    # wdLineStyleSingle
    $t.Borders.InsideLineWidth = 2    # wdLineWidth025pt -> w:sz="4"
    $t.Borders.InsideColor = 0        # wdColorBlack -> w:color="000000"
}

Write-Output "Applied borders to $($d.Tables.Count) table(s)"
